{"js": "// Office.js (Word JavaScript API) edit script.\n// Updates the job-title lines in the PROFESSIONAL EXPERIENCE section of the\n// resume to reflect the new employment history, per the commit diff.\n//\n// Exact, whole-line replacements keyed off the ORIGINAL text so the script\n// is resilient to paragraph index drift; each paragraph whose text matches\n// one of the \"before\" strings gets its text replaced with the \"after\"\n// string via a Replace-mode insertText on the paragraph range (this keeps\n// the paragraph's own formatting / style, e.g. Heading3, intact).\n\nconst replacements = [\n  {\n    before: \"PARTNER & SENIOR DATA ARCHITECT - Siege Analytics, Washington, DC | January 2014 \\u2013 Present\",\n    after: \"PARTNER - Siege Analytics, Washington, DC | January 2014 \\u2013 Present\"\n  },\n  {\n    before: \"PRINCIPAL TECHNICAL ARCHITECT - Clarity and Rigour, Washington, DC | 2012 \\u2013 2014\",\n    after: \"DATA PRODUCTS MANAGER - Helm/Murmuration, Washington, DC | 2012 \\u2013 2014\"\n  },\n  {\n    before: \"DIRECTOR OF TECHNOLOGY - Helm, Washington, DC | 2010 \\u2013 2012\",\n    after: \"SOFTWARE ENGINEER - Mautinoa Technologies, Washington, DC | 2010 \\u2013 2012\"\n  },\n  {\n    before: \"SENIOR TECHNICAL ANALYST - GSD&M, Austin, TX | 2008 \\u2013 2010\",\n    after: \"SENIOR ANALYST - Myers Research, Washington, DC | 2008 \\u2013 2010\"\n  },\n  {\n    before: \"TECHNICAL COORDINATOR - Progressive Change Campaign Committee, Washington, DC | 2006 \\u2013 2008\",\n    after: \"RESEARCH DIRECTOR - Progressive Change Campaign Committee, Washington, DC | 2006 \\u2013 2008\"\n  },\n  {\n    before: \"TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 \\u2013 2004\",\n    after: \"INTERIM TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 \\u2013 2004\"\n  },\n  {\n    before: \"TECHNICAL COORDINATOR - The Feldman Group, Washington, DC | 2000 \\u2013 2001\",\n    after: \"FIELD DIRECTOR - The Feldman Group, Washington, DC | 2000 \\u2013 2001\"\n  }\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const para of paragraphs.items) {\n  const text = para.text;\n  const match = replacements.find((r) => text === r.before);\n  if (match) {\n    para.getRange().insertText(match.after, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Updates the job-title lines in the PROFESSIONAL EXPERIENCE section of the\n# resume to reflect the new employment history, per the commit diff.\n#\n# Uses Find/Replace (Range.Find.Execute) against the whole document body for\n# each exact \"before\" line, matched case-sensitively on the WHOLE string\n# (MatchWholeWord not applicable to a full-phrase match) so lines that merely\n# share a leading job-title token (e.g. the two different \"TECHNICAL\n# COORDINATOR - ...\" entries, or \"SOFTWARE ENGINEER - ...\") are not confused\n# with one another - only the single paragraph whose text equals the \"before\"\n# value is touched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Before = \"PARTNER & SENIOR DATA ARCHITECT - Siege Analytics, Washington, DC | January 2014 \u2013 Present\"; After = \"PARTNER - Siege Analytics, Washington, DC | January 2014 \u2013 Present\" },\n    @{ Before = \"PRINCIPAL TECHNICAL ARCHITECT - Clarity and Rigour, Washington, DC | 2012 \u2013 2014\"; After = \"DATA PRODUCTS MANAGER - Helm/Murmuration, Washington, DC | 2012 \u2013 2014\" },\n    @{ Before = \"DIRECTOR OF TECHNOLOGY - Helm, Washington, DC | 2010 \u2013 2012\"; After = \"SOFTWARE ENGINEER - Mautinoa Technologies, Washington, DC | 2010 \u2013 2012\" },\n    @{ Before = \"SENIOR TECHNICAL ANALYST - GSD&M, Austin, TX | 2008 \u2013 2010\"; After = \"SENIOR ANALYST - Myers Research, Washington, DC | 2008 \u2013 2010\" },\n    @{ Before = \"TECHNICAL COORDINATOR - Progressive Change Campaign Committee, Washington, DC | 2006 \u2013 2008\"; After = \"RESEARCH DIRECTOR - Progressive Change Campaign Committee, Washington, DC | 2006 \u2013 2008\" },\n    @{ Before = \"TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 \u2013 2004\"; After = \"INTERIM TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 \u2013 2004\" },\n    @{ Before = \"TECHNICAL COORDINATOR - The Feldman Group, Washington, DC | 2000 \u2013 2001\"; After = \"FIELD DIRECTOR - The Feldman Group, Washington, DC | 2000 \u2013 2001\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #          MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    $find.Execute($r.Before, $true, $false, $false, $false, $false, $true, 1, $false, $r.After, 2) | Out-Null\n}\n"}
